# bug fix in Eduati data files
#
# The HT29_noCTRL_meas workbook had a stray block of leftover rows (45-87,
# single column A "index" values with no real data) tacked on to the bottom
# of Sheet1, left over from a previous version of the sheet that had more
# conditions. The real data only goes down to row 44. This also leaves the
# workbook with a cleaned-up view state (Sheet1 selected/active instead of
# Sheet3, and the stale selection box on Sheet1 reset to a single cell).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: drop the leftover rows 45:87 (entire rows, so everything below
#     shifts up and the used range/dimension shrinks back to A1:N44) ---
$ws1.Range("A45:A87").EntireRow.Delete()

# --- Reset Sheet1's view: scroll/select down near the old bottom edge,
#     landing on D50 (an empty cell below the real data) ---
$ws1.Range("D50").Select()

# --- Make Sheet1 the active sheet/tab again (it was Sheet3 before) ---
$ws1.Activate()
